$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $e) {
    if ($null -ne $d) { $ws.Range("D$row").Value = $d }
    if ($null -ne $e) { $ws.Range("E$row").Value = $e }
}

# Row 2 - Bitcoin
Set-Row 2 "42.056.05" "  -1.41%  "

# Row 3 - Ethereum
Set-Row 3 "2.303.70" "  -2.27%  "

# Row 4 - TetherUSD
Set-Row 4 $null "  -0.13%  "

# Row 5 - BNB
Set-Row 5 "318.82" "  +0.33%  "

# Row 6 - Solana
Set-Row 6 "104.94" "  -3.08%  "

# Row 7 - XRP
Set-Row 7 "0.631" "  -0.88%  "

# Row 9 - Cardano
Set-Row 9 $null "  -1.49%  "

# Row 10 - Avalanche
Set-Row 10 "39.71" "  -4.88%  "

# Row 11 - Dogecoin
Set-Row 11 "0.0912" "  -1.63%  "

# Row 12 - Polkadot
Set-Row 12 $null "  -0.75%  "

# Row 13 - TRON
Set-Row 13 $null "  +0.28%  "

# Row 14 - Polygon
Set-Row 14 "0.981" "  -1.75%  "

# Row 15 - Chainlink
Set-Row 15 "15.44" "  -3.32%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-Row 16 "2.652.11" "  -2.29%  "

# Row 17 - WrappedEther
Set-Row 17 "2.305.32" "  -2.58%  "

# Row 18 - WrappedBTC
Set-Row 18 "42.043.92" "  -1.49%  "

# Row 19 - Uniswap
Set-Row 19 "7.74" "  +0.10%  "

# Row 20 - ShibaInu
Set-Row 20 $null "  -0.58%  "

# Row 21 - BitcoinCash
Set-Row 21 "286.74" "  +11.43%  "

# Row 22 - Litecoin
Set-Row 22 "73.64" $null

# Row 23 - PancakeSwap
Set-Row 23 "3.62" "  +1.01%  "

# Rows 24/25 swap: InternetComputer(DFINITY) <-> ImmutableX
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  -1.60%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "10.06"
$ws.Range("E25").Value = "  +6.76%  "

# Row 26 - Dai
Set-Row 26 "1.01" "  +0.62%  "

# Row 27 - Cosmos
Set-Row 27 "10.96" "  -4.25%  "

# Row 28 - EthereumClassic
Set-Row 28 "23.44" "  +2.70%  "

# Row 29 - Toncoin
Set-Row 29 "2.27" "  +1.66%  "

# Rows 30/31 swap: Monero <-> InjectiveProtocol
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "35.52"
$ws.Range("E30").Value = "  -3.48%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "164.33"
$ws.Range("E31").Value = "  -5.93%  "

# Row 32 - Hedera
Set-Row 32 "0.0887" "  -0.45%  "

# Row 33 - WEMIXToken
Set-Row 33 $null "  +1.13%  "

# Row 34 - Filecoin
Set-Row 34 "5.89" "  -3.30%  "

# Row 35 - Stellar
Set-Row 35 $null "  +0.77%  "

# Row 36 - Kaspa
Set-Row 36 $null "  -8.83%  "

# Row 37 - RenderToken
Set-Row 37 "4.64" "  +0.51%  "

# Row 38 - LidoDAOToken
Set-Row 38 "2.94" "  +9.94%  "

# Row 39 - VeChain
Set-Row 39 "0.0352" "  -2.88%  "

# Row 40 - NEARProtocol
Set-Row 40 $null "  -5.76%  "

# Row 41 - BitcoinSV
Set-Row 41 "101.79" "  +19.85%  "

# Row 42 - ARBITRUM
Set-Row 42 $null "  +1.51%  "

# Row 43 - MultiversX
Set-Row 43 "71.04" "  -1.10%  "

# Row 44 - Algorand
Set-Row 44 $null "  -4.34%  "

# Row 45 - FirstDigitalUSD
Set-Row 45 $null "  -0.03%  "

# Row 46 - Aave
Set-Row 46 "117.14" "  +2.85%  "

# Row 47 - Celestia
Set-Row 47 $null "  +1.25%  "

# Row 48 - FraxShare
Set-Row 48 "9.15" "  +0.10%  "

# Row 49 - ordi
Set-Row 49 "77.96" "  +4.93%  "

# Row 50 - THORChain
Set-Row 50 "5.36" "  -2.54%  "

# Row 51 - TrustWalletToken
Set-Row 51 "1.29" "  +0.96%  "
